$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "study title test1"
$ws.Range("B2").Value = "sample1"
$ws.Range("C2").Value = "Human"
$ws.Range("E2").Value = "Plasma"
$ws.Range("F2").Value = "Something"
$ws.Range("G2").Value = "SK"
$ws.Range("H2").Value = "SAF"
$ws.Range("J2").Value = "'123"
$ws.Range("L2").Value = "'45678"
$ws.Range("M2").Value = "human id"
$ws.Range("N2").Value = "human1"
$ws.Range("O2").Value = "Consent form is here..."
$ws.Range("P2").Value = "THIS-IS-A-CODE"
$ws.Range("Q2").Value = "'1"
$ws.Range("R2").Value = "test"
$ws.Range("I2").Value = "G134"
$ws.Range("D2").Value = "Eppendorf"
$ws.Range("K2").Value = 1
